# Lesson 25 (slide 4) - fix the interrupt example code so the timer
# configuration bit names match the TI MSP430 "_L"-suffixed style used
# elsewhere in the deck:
#   TASSEL1      -> TASSEL_2
#   ID1|ID0      -> ID_3            (also tidy the trailing whitespace to a tab)
#   MC1          -> MC_1
#
# The target shape is the "Content Placeholder 2" code listing on slide 4.
# Each change lives in its own paragraph, so we locate the paragraph by its
# (stable) current text, then surgically replace just the changed substring
# via TextRange.Characters(start, length) so the surrounding run formatting
# (black for code / green for the "//" comment) is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

function Replace-InParagraph($TextRange, $ParaIndex, $OldSubstring, $NewSubstring) {
    $para = $TextRange.Paragraphs($ParaIndex, 1)
    $text = $para.Text
    $idx = $text.IndexOf($OldSubstring)
    if ($idx -lt 0) {
        throw "Could not find '$OldSubstring' in paragraph $ParaIndex (text: '$text')"
    }
    $chars = $para.Characters($idx + 1, $OldSubstring.Length)
    $chars.Text = $NewSubstring
}

# Paragraph 9: "    TA0CTL |= TASSEL1;           // configure for SMCLK - what's the frequency (roughly)?"
Replace-InParagraph $tr 9 "TASSEL1" "TASSEL_2"

# Paragraph 10: "    TA0CTL |= ID1|ID0;           // divide clock by 8 - what's the frequency of interrupt?"
Replace-InParagraph $tr 10 "ID1|ID0;           " "ID_3;         `t  "

# Paragraph 12: "    TA0CTL |= MC1;               // set count mode to continuous"
Replace-InParagraph $tr 12 "MC1" "MC_1"
